$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "WalgreensDataSheet"

# Set cell values (order chosen to match shared-string insertion order: Item, Key, Value, Eye Drops)
$ws.Range("A2").Value = "Item"
$ws.Range("A1").Value = "Key"
$ws.Range("B1").Value = "Value"
$ws.Range("B2").Value = "Eye Drops"

# Apply a thin border around every used cell A1:B5
$allRange = $ws.Range("A1:B5")
$allRange.Borders.LineStyle = 1

# Header formatting: bold font + yellow fill, applied on a single cell first
# then propagated via copy/paste-special so no stray intermediate cell style
# gets created (each direct Range.Font / Range.Interior write on a multi-cell
# range creates its own style entry).
$a1 = $ws.Range("A1")
$a1.Font.Bold = $true
$a1.Interior.Color = 65535

$a1.Copy()
$ws.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the author's final selection
$ws.Range("B2").Select() | Out-Null
